# Applies crypto price/volume/coin-list updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'301.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.00%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'32.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.51%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.104"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.34%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.07708"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.75%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'2.138"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.07%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'7.857"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.61%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9251"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.30%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1766"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.45%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.07963"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.54%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08466"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-4.59%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03066"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.03%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09994"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.25%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001525"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.99%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005849"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-2.81%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007498"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2,116.77%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.471"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.10%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.780"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.70%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'2.154"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-4.27%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'0.3347"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.66%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.1314"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.05%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'4.270"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.40%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.1969"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'9.91%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.04524"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-1.61%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("E25").Value = "'-1.51%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.004851"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'8.25%"
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.0001248"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.10%"
$ws.Range("E27").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.01717"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.92%"
$ws.Range("E39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.04702"
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007515"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.62%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1362"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.22%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.92%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.01059"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.62%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00006233"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-0.19%"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.11%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'1.042"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'49.51%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.002996"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'19.38%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.00002097"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.11%"
$ws.Range("E49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.0001997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.11%"
$ws.Range("E50").Style = "Normal"

Write-Host "Applied crypto symbol list update"
